$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.519.14'
$ws.Range('E2').Value = '  -0.51%  '
$ws.Range('D3').Value = '2.948.87'
$ws.Range('E3').Value = '  -1.97%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.94'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.04'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.36%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  -0.30%  '
$ws.Range('D9').Value = '2.945.56'
$ws.Range('E9').Value = '  -1.96%  '
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('E11').Value = '  -4.21%  '
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000244'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -2.74%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.84'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +0.35%  '
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '65.539.39'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').Value = '3.437.31'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.05'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +0.98%  '
$ws.Range('D19').Value = '2.945.47'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.74'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +12.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '446.63'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.695'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.29'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -1.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.13'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.24'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -2.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.26'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.90%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.03'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -6.18%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.49'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +6.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.04'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.60'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('E33').Value = '  +3.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.22'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.972'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -1.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.72'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '46.21'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +5.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.08'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.304'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -1.26%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.97'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -8.12%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.82'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -5.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.51'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '382.77'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = '2.679.86'
$ws.Range('E47').Value = '  -4.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '133.46'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.87'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.17'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +1.27%  '
